$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 74 values (date in A74 stays "01-01-2021")
$ws.Range("B74").Value = 611
$ws.Range("C74").Value = -210
$ws.Range("D74").Value = 10214
$ws.Range("E74").Value = -1299
$ws.Range("F74").Value = -1165
$ws.Range("G74").Value = 8151
$ws.Range("H74").Value = 4657
$ws.Range("I74").Value = 3494
$ws.Range("J74").Value = 3340
$ws.Range("K74").Value = 154
$ws.Range("L74").Value = 674
$ws.Range("M74").Value = -281
$ws.Range("N74").Value = -802

# Add new row 75
# (build the date-like label as plain text without letting Excel's
# autodetect convert it to a real date/serial number or add a style)
$ws.Range("A75").Formula = '="01-04-2021"'
$ws.Range("A75").Copy()
$ws.Range("A75").PasteSpecial(-4163)
$ws.Range("B75").Value = 629
$ws.Range("C75").Value = -43
$ws.Range("D75").Value = 10250
$ws.Range("E75").Value = -1327
$ws.Range("F75").Value = -4352
$ws.Range("G75").Value = 5158
$ws.Range("H75").Value = 5186
$ws.Range("I75").Value = -28
$ws.Range("J75").Value = 3309
$ws.Range("K75").Value = -3336
$ws.Range("L75").Value = 982
$ws.Range("M75").Value = -708
$ws.Range("N75").Value = -5026
